$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.519.74'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '3.116.69'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.117.18'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.448'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.398'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.37%  '
$ws.Range('D13').Value = '3.654.51'
$ws.Range('E14').Value = '  +2.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '57.657.70'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = '3.114.78'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '348.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.37'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.72%  '
$ws.Range('E26').Value = '  -1.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.168'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').Value = '0.0₃0913'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('E30').Value = '  +3.57%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('E35').Value = '  -1.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.97'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.15'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.20'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.76%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0666'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.70%  '
$ws.Range('E43').Value = '  +7.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.701'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('D45').Value = '3.157.40'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '36.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0269'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.71%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.344.15'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.959'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('E51').Value = '  +0.47%  '
